$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 302, shifting existing rows 302-396 down to 304-398
$ws.Rows("302:303").Insert()

# Fill in the two new rows with the new weekly record (date 44809), same fixed
# columns as the rest of the Coliflor / Feria Lagunitas de Puerto Montt block.
$ws.Cells.Item(302, 1).Value = 4
$ws.Cells.Item(302, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(302, 3).Value = "Los Lagos"
$ws.Cells.Item(302, 4).Value = 44809
$ws.Cells.Item(302, 5).Value = 10
$ws.Cells.Item(302, 6).Value = 100112008
$ws.Cells.Item(302, 7).Value = "Coliflor"
$ws.Cells.Item(302, 8).Value = "Sin especificar"
$ws.Cells.Item(302, 9).Value = "Primera"
$ws.Cells.Item(302, 10).Value = 250
$ws.Cells.Item(302, 11).Value = 2000
$ws.Cells.Item(302, 12).Value = 2000
$ws.Cells.Item(302, 13).Value = 2000
$ws.Cells.Item(302, 14).Value = "`$/unidad"
$ws.Cells.Item(302, 15).Value = "Región del Maule"
$ws.Cells.Item(302, 16).Value = 2000
$ws.Cells.Item(302, 17).Value = 1
$ws.Cells.Item(302, 18).Value = "Hortaliza"

$ws.Cells.Item(303, 1).Value = 4
$ws.Cells.Item(303, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(303, 3).Value = "Los Lagos"
$ws.Cells.Item(303, 4).Value = 44809
$ws.Cells.Item(303, 5).Value = 10
$ws.Cells.Item(303, 6).Value = 100112008
$ws.Cells.Item(303, 7).Value = "Coliflor"
$ws.Cells.Item(303, 8).Value = "Sin especificar"
$ws.Cells.Item(303, 9).Value = "Segunda"
$ws.Cells.Item(303, 10).Value = 250
$ws.Cells.Item(303, 11).Value = 1800
$ws.Cells.Item(303, 12).Value = 1800
$ws.Cells.Item(303, 13).Value = 1800
$ws.Cells.Item(303, 14).Value = "`$/unidad"
$ws.Cells.Item(303, 15).Value = "Región del Maule"
$ws.Cells.Item(303, 16).Value = 1800
$ws.Cells.Item(303, 17).Value = 1
$ws.Cells.Item(303, 18).Value = "Hortaliza"

# Apply the date/time number format used by the rest of the Fecha (D) column
$ws.Range("D302:D303").NumberFormat = "YYYY-MM-DD HH:MM:SS"
